# Revert changes for SampleData:
#  - Clear the "source" column (C) values for rows 246-263 (these cells are
#    removed entirely from the sheet, matching how Excel drops an empty,
#    unformatted cell when its contents are cleared).
#  - Update the active sheet's selection to the single cell I258.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C246:C263").ClearContents()

$ws.Range("I258").Select()
